# Input data check for sensor errors
# Created from Panos

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Parameters")

# --- New header / data in columns Q:R -------------------------------------
# Values are written in the same order the original author entered them so
# the shared-string table append order matches the target file exactly.
$ws.Range("B1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value = "Parameters that have to be check before HC "

$ws.Range("B1").Copy()
$ws.Range("Q14").PasteSpecial(-4122)
$ws.Range("Q14").Value = "Parameters that have to be checked after HC"

$ws.Range("Q2").Value = "AF pressure drop"
$ws.Range("Q3").Value = "Pscav - Pexh"
$ws.Range("Q4").Value = "FPI vs Power"

# --- New row 28/29 data in column F ----------------------------------------
$ws.Range("F14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = "Parameters that should be compared in absolut numbers and not in percentage"
$ws.Rows(28).RowHeight = 15.75

$ws.Range("Q17").Value = "Texh"

$ws.Range("Q5").Value = "Energy equilibrium on TC shaft (C - T)"

$ws.Range("R15").Value = "usually corrected to ISO using Tamb and TcwIN"
$ws.Range("R16").Value = "usually corrected to ISO using Tamb and TcwIN"
$ws.Range("R17").Value = "usually corrected to ISO using Tamb and TcwIN"
$ws.Range("R18").Value = "usually corrected to ISO using Tamb and TcwIN"

$ws.Range("Q15").Value = "Pmax"
$ws.Range("Q16").Value = "Pcomp"
$ws.Range("Q18").Value = "Pscav"

$ws.Range("F29").Value = "Tscav"

# --- Column widths -----------------------------------------------------
# (ColumnWidth rounds internally to 1/6-character steps, so the inputs below
# are chosen to land as close as possible on the target stored widths.)
$ws.Columns.Item(3).ColumnWidth = 13
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(17).ColumnWidth = 45.166666666666664
$ws.Columns.Item(18).ColumnWidth = 49.833333333333336

# --- View state: make "All Parameters" the active / selected tab -----------
$ws.Activate()
$ws.Range("R15").Select()

# --- "Component diagnosis" column E width (manual resize, no longer bestFit)
$ws3 = $wb.Worksheets.Item("Component diagnosis")
$ws3.Columns.Item(5).ColumnWidth = 32.666666666666664

# Re-activate "All Parameters" so it stays the selected/active sheet tab.
$ws.Activate()
